# Run for CCP 1vs1 IM+DF=0.8 OK
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil2")
$ws.Activate()

# --- First block (rows 4-9, "n=1"): fill column D with the new values ---
$ws.Range("D4").Value = 0.03475217
$ws.Range("D5").Value = 0.0173505288021
$ws.Range("D6").Value = -0.0472881790235
$ws.Range("D7").Value = 0.0484820911641
$ws.Range("D8").Value = 0.0558960777764
$ws.Range("D9").Value = 0.0592033378977
# D10:D12 are a shared formula (=SUM(D$4:D$6, D7/D8/D9)); they recalc on their own.

# --- Second block (rows 16-21, "n=8"): fill column D with the new values ---
$ws.Range("D16").Value = 0.13220419
$ws.Range("D17").Value = 0.112277716654
$ws.Range("D18").Value = -0.0714823681627

# D19 and D21 swap their cell styles (D19: 49->50, D21: 50->49). Copy the
# formats from the rows that already carry those styles (D9 has style 50,
# D5 has style 49), then (re)apply the values.
$ws.Cells.Item(9, 4).Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4122)
$ws.Cells.Item(5, 4).Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D19").Value = 0.0208765988749
$ws.Range("D20").Value = 0.0189590500981
$ws.Range("D21").Value = 0.0181729894731

# E19 becomes completely empty (no style, no content) - the <c> disappears.
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(19, 5).ClearContents()

# --- Update the SUM formulas in D22:D24 to reference the rotated rows ---
$ws.Range("D22").Formula = "=SUM(D`$4:D`$6, D20)"
$ws.Range("D23").Formula = "=SUM(D`$4:D`$6, D21)"
$ws.Range("D24").Formula = "=SUM(D`$4:D`$6, D19)"

# --- View state: scroll so row 7 is at the top, select D7 ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D7").Select()
